# Add 2022-Q3 data:
#  - Duplicate the "2022-Q2" sheet (keeps header/style layout), rename it
#    to "2022-Q3", place it right before "2022-Q2", trim it down to the
#    Q3 row count, and overwrite the fund rows with the Q3 figures.
#  - Insert a new summary row on "总计" for 2022-Q3, above the existing
#    2022-Q2 / 2022-Q1 rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet by duplicating "2022-Q2" (so the
#    header text / bold-bordered styling comes along for free) and
#    dropping it in right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# 2022-Q2 has 14 fund rows (rows 2-15); 2022-Q3 only has 12 (rows 2-13),
# so drop the two extra trailing rows that got copied along.
$q3.Rows("14:15").Delete()

$q3data = @(
    @("0","014294","南方北交所精选两年定开混合","4.26","75.23","6.26","0.2667","3"),
    @("1","014271","大成北交所两年定开混合A","3.45","65.31","6.12","0.2111","6"),
    @("2","014283","华夏北交所创新中小企业精选两年定开混合","3.47","71.59","4.79","0.1662","7"),
    @("3","014273","广发北交所精选两年定开混合A","3.37","64.25","4.77","0.1607","6"),
    @("4","014277","万家北交所慧选两年定期开放混合A","3.56","93.97","4.09","0.1456","10"),
    @("5","014279","汇添富北交所创新精选两年定开混合A","3.20","93.27","3.60","0.1152","10"),
    @("6","014275","易方达北交所精选两年定开混合A","3.58","61.75","3.17","0.1135","4"),
    @("7","014272","大成北交所两年定开混合C","0.82","65.31","6.12","0.0502","6"),
    @("8","014274","广发北交所精选两年定开混合C","0.85","64.25","4.77","0.0405","6"),
    @("9","014276","易方达北交所精选两年定开混合C","0.92","61.75","3.17","0.0292","4"),
    @("10","014278","万家北交所慧选两年定期开放混合C","0.49","93.97","4.09","0.0200","10"),
    @("11","014280","汇添富北交所创新精选两年定开混合C","0.51","93.27","3.60","0.0184","10")
)

# Columns B-G hold fund codes / names / ratios that look numeric (leading
# zeros, decimals) but must stay text, matching the source data's
# inlineStr typing. Force Text format before writing, then strip the
# formatting override back off (the underlying rows already inherited
# plain/no-style formatting from the 2022-Q2 copy) so only the value
# type - not an extra style - changes.
$q3.Range("B2:G13").NumberFormat = "@"
for ($i = 0; $i -lt $q3data.Length; $i++) {
    $row = $q3data[$i]
    $r = 2 + $i
    $q3.Range("A$r").Value = [double]$row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").Value = $row[6]
    $q3.Range("H$r").Value = [double]$row[7]
}
$q3.Range("B2:G13").ClearFormats()

# ---------------------------------------------------------------------
# 2) Update "总计": insert a new row 2 for 2022-Q3 above the existing
#    2022-Q2 / 2022-Q1 rows (which shift down to rows 3 / 4).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# The insert copies row-1's (header) formatting onto the new blank row;
# strip that off the B:D cells (no style in the target), and pull A's
# index-column style from the row below instead.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 1.34

# The A column is a 0-based running index; renumber the rows that shifted
# down (old index 0/1 -> new index 1/2).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# Copying "2022-Q2" left the new "2022-Q3" sheet active/selected. Put the
# active tab back on "2022-Q1" (the sheet that was selected originally)
# so the per-sheet view state is unaffected by this edit.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()
